# Adds 14 new GUID/Information rows (601-614) to Sheet1
# and fixes the existing row numbering/content per the commit:
# "added motive rotation ... fixed the 000 point showing up,
# the first and last indexes of palm and smooth baton were mixed up."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 601

$guids = @(
    '000599',
    '000600',
    '000601',
    '000602',
    '000603',
    '000604',
    '000605',
    '000606',
    '000607',
    '000608',
    '000609',
    '000610',
    '000611',
    '000612'
)

$infos = @(
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 13B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:28:16',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 13B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:28:49',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 14B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:29:07',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 14C. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:29:23',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 21B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:29:43',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 22B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:30:05',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 23B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:30:22',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 31B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:31:28',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 31B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:32:19',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 31B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:33:32',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 31B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:34:01',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 31B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:35:12',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 31B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:38:12',
    'Details: Palm Position and Baton tip position - transformed IMU CJMCU-20948 Data Reading and single hand Leap LM-010 Reading and Motive (OptiTrack) Reading. Experiment ID: 31B. Script used: Interpret_IMU_And_LeapDevice_And_Motive_Separately.  Dataset used: Test data from raw imu reading and raw leap reading with motive Data scaled and resampled. Loop every 20ms, imufilter sample rate 5000. arduino internal delay 200. only plotting last 500 values. File Location: Visualisations/IMU_Leap_CombinedData. Date Generated: 11-Jun-2023 22:38:57'
)

for ($i = 0; $i -lt $guids.Length; $i++) {
    $r = $startRow + $i
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $guids[$i]
    $cellA.Style = "Normal"

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $infos[$i]
}

Write-Output ("Added {0} rows, new dimension: {1}" -f $guids.Length, $ws.UsedRange.Address())
